$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "MATRICULA"
$ws.Range("B1").Value = "NOMBRE"
$ws.Range("C1").Value = "CURP"

# --- Data row (row 2) ---
$ws.Range("A2").Value = "189625"
$ws.Range("B2").Value = "Carlos Segoviano"
$ws.Range("C2").Value = "SACC000117HSRLLRA5"

# --- Column widths (values chosen so the engine's pixel-quantized stored
#     width lands on the closest possible match to the template's
#     19.28515625 / 26.140625 / 28.42578125 character widths) ---
$ws.Columns("A").ColumnWidth = 18.416666666666668
$ws.Columns("B").ColumnWidth = 25.25
$ws.Columns("C").ColumnWidth = 27.583333333333332

# --- Formatting: Arial + centered for the whole used block ---
$dataRange = $ws.Range("A1:D22")
$dataRange.Font.Name = "Arial"
$dataRange.HorizontalAlignment = -4108   # xlCenter

# --- C3 keeps an underline (like the template's old placeholder cell) ---
$ws.Range("C3").Font.Underline = $true

# --- Selection matches the saved file ---
$ws.Range("C2").Select()

Write-Host "done"
